$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 1; $r -le 3; $r++) {
    $full = $ws.Cells.Item($r, 2).Value2
    $parts = $full -split "`n"
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $ws.Cells.Item($r, 2 + $i).Value = $parts[$i]
    }
}
